$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 1: in the H5Ocopy code sample, Word's "_GoBack" bookmark had split
# the text "status = H5Ocopy(file1_id, ...)" into two runs around the
# cursor's last-edit position. Re-typing/replacing the full line merges
# it back into a single run and drops that now-stale bookmark.
# -----------------------------------------------------------------------
$rngMerge = $d.Content
$rngMerge.Start = 0
$foundMerge = $rngMerge.Find
$foundMerge.ClearFormatting()
$foundMerge.Execute( `
    "status = H5Ocopy(file1_id, src_name, file2_id, dst_name, ocpypl_id, H5P_DEFAULT);", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "status = H5Ocopy(file1_id, src_name, file2_id, dst_name, ocpypl_id, H5P_DEFAULT);", 2) | Out-Null

# -----------------------------------------------------------------------
# Edit 2: in the "Function Summary" table, the Fortran-equivalent cell
# for H5Ocopy said "(none)"; fill in the real Fortran wrapper name.
# Locate the first "(none)" that follows the "H5Ocopy" table entry
# (skipping all the earlier prose mentions of H5Ocopy).
# -----------------------------------------------------------------------
$rng = $d.Range(10900, $d.Content.End)
$found = $rng.Find
$found.ClearFormatting()
$found.Text = "H5Ocopy"
$found.Forward = $true
$found.Wrap = 0
$found.Execute() | Out-Null

$rng2 = $d.Range($rng.End, $rng.End + 40)
$found2 = $rng2.Find
$found2.ClearFormatting()
$found2.Text = "(none)"
$found2.Forward = $true
$found2.Wrap = 0
$found2.Execute() | Out-Null

# Replace the text, keeping one extra placeholder character at the end so
# that the bookmark we add next lands on an ordinary text position rather
# than exactly on the paragraph-end mark.
$rng2.Text = "H5ocopy_fZ"

$bm = $d.Range($rng2.End - 1, $rng2.End - 1)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null

# Strip the placeholder character again; the bookmark stays put, now
# sitting right after "H5ocopy_f" at the paragraph end - matching what
# Word leaves behind after you type the replacement text there.
$delRange = $d.Range($rng2.End - 1, $rng2.End)
$delRange.Text = ""
